# Auto-generated Excel COM-interop edit script
# Applies numeric value corrections to market-profit tables across all 8 item-sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the source diff.

$wb = $excel.ActiveWorkbook

# --- Cell value updates (Sheet, Cell, NewValue) ---
$updates = @(
    @{Sheet="ALC"; Cell="H108"; Value=31620}
    @{Sheet="ALC"; Cell="J108"; Value=31620}
    @{Sheet="ALC"; Cell="L108"; Value=31620}
    @{Sheet="ALC"; Cell="N108"; Value=-39300}
    @{Sheet="ALC"; Cell="H116"; Value=6112.3335}
    @{Sheet="ALC"; Cell="I116"; Value=3800}
    @{Sheet="ALC"; Cell="K116"; Value=3800}
    @{Sheet="ALC"; Cell="M116"; Value=-358}
    @{Sheet="ALC"; Cell="H117"; Value=48538}
    @{Sheet="ALC"; Cell="J117"; Value=48538}
    @{Sheet="ALC"; Cell="L117"; Value=48538}
    @{Sheet="ALC"; Cell="N117"; Value=-57716}
    @{Sheet="ALC"; Cell="H128"; Value=38311}
    @{Sheet="ALC"; Cell="J128"; Value=38311}
    @{Sheet="ALC"; Cell="L128"; Value=38311}
    @{Sheet="ALC"; Cell="N128"; Value=-48271}
    @{Sheet="ALC"; Cell="H133"; Value=38279.09}
    @{Sheet="ALC"; Cell="J133"; Value=38279.09}
    @{Sheet="ALC"; Cell="L133"; Value=38279.09}
    @{Sheet="ALC"; Cell="N133"; Value=-48399.09}
    @{Sheet="ALC"; Cell="H135"; Value=11364512}
    @{Sheet="ALC"; Cell="I135"; Value=716}
    @{Sheet="ALC"; Cell="J135"; Value=71430290}
    @{Sheet="ALC"; Cell="K135"; Value=6444}
    @{Sheet="ALC"; Cell="L135"; Value=642872610}
    @{Sheet="ALC"; Cell="M135"; Value=-3909}
    @{Sheet="ALC"; Cell="N135"; Value=-642877680}
    @{Sheet="ALC"; Cell="H136"; Value=36131.785}
    @{Sheet="ALC"; Cell="J136"; Value=36131.785}
    @{Sheet="ALC"; Cell="L136"; Value=36131.785}
    @{Sheet="ALC"; Cell="N136"; Value=-46331.785}
    @{Sheet="ALC"; Cell="H139"; Value=45413.266}
    @{Sheet="ALC"; Cell="J139"; Value=45413.266}
    @{Sheet="ALC"; Cell="L139"; Value=45413.266}
    @{Sheet="ALC"; Cell="N139"; Value=-55693.266}
    @{Sheet="ARM"; Cell="H32"; Value=28149.344}
    @{Sheet="ARM"; Cell="I32"; Value=27775.34}
    @{Sheet="ARM"; Cell="K32"; Value=27775.34}
    @{Sheet="ARM"; Cell="M32"; Value=-27488.34}
    @{Sheet="ARM"; Cell="H61"; Value=2682.1277}
    @{Sheet="ARM"; Cell="I61"; Value=1511.5862}
    @{Sheet="ARM"; Cell="K61"; Value=1511.5862}
    @{Sheet="ARM"; Cell="M61"; Value=-1299.5862}
    @{Sheet="ARM"; Cell="H123"; Value=0}
    @{Sheet="ARM"; Cell="J123"; Value=0}
    @{Sheet="ARM"; Cell="L123"; Value=0}
    @{Sheet="ARM"; Cell="H125"; Value=50715}
    @{Sheet="ARM"; Cell="J125"; Value=50715}
    @{Sheet="ARM"; Cell="L125"; Value=50715}
    @{Sheet="ARM"; Cell="N125"; Value=-60555}
    @{Sheet="ARM"; Cell="H130"; Value=39852.668}
    @{Sheet="ARM"; Cell="J130"; Value=39852.668}
    @{Sheet="ARM"; Cell="L130"; Value=39852.668}
    @{Sheet="ARM"; Cell="N130"; Value=-49892.668}
    @{Sheet="ARM"; Cell="H133"; Value=42745.25}
    @{Sheet="ARM"; Cell="J133"; Value=42745.25}
    @{Sheet="ARM"; Cell="L133"; Value=42745.25}
    @{Sheet="ARM"; Cell="N133"; Value=-47805.25}
    @{Sheet="ARM"; Cell="H134"; Value=52224}
    @{Sheet="ARM"; Cell="J134"; Value=52224}
    @{Sheet="ARM"; Cell="L134"; Value=52224}
    @{Sheet="ARM"; Cell="N134"; Value=-62364}
    @{Sheet="ARM"; Cell="H136"; Value=2682.1277}
    @{Sheet="ARM"; Cell="I136"; Value=1511.5862}
    @{Sheet="ARM"; Cell="K136"; Value=4534.7586}
    @{Sheet="ARM"; Cell="M136"; Value=-1984.7586}
    @{Sheet="BSM"; Cell="H81"; Value=8473.333000000001}
    @{Sheet="BSM"; Cell="J81"; Value=8473.333000000001}
    @{Sheet="BSM"; Cell="L81"; Value=8473.333000000001}
    @{Sheet="BSM"; Cell="N81"; Value=-10595.333}
    @{Sheet="BSM"; Cell="H84"; Value=8473.333000000001}
    @{Sheet="BSM"; Cell="J84"; Value=8473.333000000001}
    @{Sheet="BSM"; Cell="L84"; Value=25419.999}
    @{Sheet="BSM"; Cell="N84"; Value=-36027.999}
    @{Sheet="BSM"; Cell="H132"; Value=0}
    @{Sheet="BSM"; Cell="J132"; Value=0}
    @{Sheet="BSM"; Cell="L132"; Value=0}
    @{Sheet="BSM"; Cell="H134"; Value=3660.8105}
    @{Sheet="BSM"; Cell="I134"; Value=1510.0857}
    @{Sheet="BSM"; Cell="J134"; Value=4915.4}
    @{Sheet="BSM"; Cell="K134"; Value=4530.257100000001}
    @{Sheet="BSM"; Cell="L134"; Value=14746.2}
    @{Sheet="BSM"; Cell="M134"; Value=-1995.257100000001}
    @{Sheet="BSM"; Cell="N134"; Value=-19816.2}
    @{Sheet="BSM"; Cell="H135"; Value=0}
    @{Sheet="BSM"; Cell="J135"; Value=0}
    @{Sheet="BSM"; Cell="L135"; Value=0}
    @{Sheet="BSM"; Cell="H137"; Value=0}
    @{Sheet="BSM"; Cell="J137"; Value=0}
    @{Sheet="BSM"; Cell="L137"; Value=0}
    @{Sheet="CRP"; Cell="H13"; Value=350}
    @{Sheet="CRP"; Cell="I13"; Value=200}
    @{Sheet="CRP"; Cell="J13"; Value=500}
    @{Sheet="CRP"; Cell="K13"; Value=200}
    @{Sheet="CRP"; Cell="L13"; Value=500}
    @{Sheet="CRP"; Cell="M13"; Value=-61}
    @{Sheet="CRP"; Cell="N13"; Value=-778}
    @{Sheet="CRP"; Cell="H52"; Value=41999.5}
    @{Sheet="CRP"; Cell="J52"; Value=41999.5}
    @{Sheet="CRP"; Cell="L52"; Value=41999.5}
    @{Sheet="CRP"; Cell="N52"; Value=-42587.5}
    @{Sheet="CRP"; Cell="H58"; Value=1524.9744}
    @{Sheet="CRP"; Cell="I58"; Value=1315.129}
    @{Sheet="CRP"; Cell="J58"; Value=2338.125}
    @{Sheet="CRP"; Cell="K58"; Value=1315.129}
    @{Sheet="CRP"; Cell="L58"; Value=2338.125}
    @{Sheet="CRP"; Cell="M58"; Value=-1112.129}
    @{Sheet="CRP"; Cell="N58"; Value=-2744.125}
    @{Sheet="CRP"; Cell="H100"; Value=37443.332}
    @{Sheet="CRP"; Cell="J100"; Value=37443.332}
    @{Sheet="CRP"; Cell="L100"; Value=37443.332}
    @{Sheet="CRP"; Cell="N100"; Value=-39607.332}
    @{Sheet="CRP"; Cell="H135"; Value=26170}
    @{Sheet="CRP"; Cell="J135"; Value=26170}
    @{Sheet="CRP"; Cell="L135"; Value=26170}
    @{Sheet="CRP"; Cell="N135"; Value=-36310}
    @{Sheet="CRP"; Cell="H136"; Value=1524.9744}
    @{Sheet="CRP"; Cell="I136"; Value=1315.129}
    @{Sheet="CRP"; Cell="J136"; Value=2338.125}
    @{Sheet="CRP"; Cell="K136"; Value=3945.387}
    @{Sheet="CRP"; Cell="L136"; Value=7014.375}
    @{Sheet="CRP"; Cell="M136"; Value=-1395.387}
    @{Sheet="CRP"; Cell="N136"; Value=-12114.375}
    @{Sheet="CRP"; Cell="H137"; Value=31893.334}
    @{Sheet="CRP"; Cell="J137"; Value=31893.334}
    @{Sheet="CRP"; Cell="L137"; Value=31893.334}
    @{Sheet="CRP"; Cell="N137"; Value=-42093.334}
    @{Sheet="CRP"; Cell="H138"; Value=51999}
    @{Sheet="CRP"; Cell="J138"; Value=51999}
    @{Sheet="CRP"; Cell="L138"; Value=51999}
    @{Sheet="CRP"; Cell="N138"; Value=-62279}
    @{Sheet="CRP"; Cell="H141"; Value=13672.5}
    @{Sheet="CRP"; Cell="J141"; Value=15900}
    @{Sheet="CRP"; Cell="L141"; Value=15900}
    @{Sheet="CRP"; Cell="N141"; Value=-26260}
    @{Sheet="CUL"; Cell="H8"; Value=1527.0714}
    @{Sheet="CUL"; Cell="I8"; Value=1527.0714}
    @{Sheet="CUL"; Cell="K8"; Value=4581.2142}
    @{Sheet="CUL"; Cell="M8"; Value=-4442.2142}
    @{Sheet="GSM"; Cell="H110"; Value=47997}
    @{Sheet="GSM"; Cell="J110"; Value=47997}
    @{Sheet="GSM"; Cell="L110"; Value=47997}
    @{Sheet="GSM"; Cell="N110"; Value=-56177}
    @{Sheet="GSM"; Cell="H122"; Value=1800.091}
    @{Sheet="GSM"; Cell="I122"; Value=1766.6666}
    @{Sheet="GSM"; Cell="J122"; Value=1840.2}
    @{Sheet="GSM"; Cell="K122"; Value=5299.9998}
    @{Sheet="GSM"; Cell="L122"; Value=5520.6}
    @{Sheet="GSM"; Cell="M122"; Value=-2849.9998}
    @{Sheet="GSM"; Cell="N122"; Value=-10420.6}
    @{Sheet="GSM"; Cell="H130"; Value=48887.43}
    @{Sheet="GSM"; Cell="J130"; Value=48887.43}
    @{Sheet="GSM"; Cell="L130"; Value=48887.43}
    @{Sheet="GSM"; Cell="N130"; Value=-58927.43}
    @{Sheet="GSM"; Cell="H135"; Value=48000}
    @{Sheet="GSM"; Cell="J135"; Value=48000}
    @{Sheet="GSM"; Cell="L135"; Value=48000}
    @{Sheet="GSM"; Cell="N135"; Value=-58140}
    @{Sheet="LTW"; Cell="H93"; Value=1986.9375}
    @{Sheet="LTW"; Cell="J93"; Value=2129.5}
    @{Sheet="LTW"; Cell="L93"; Value=2129.5}
    @{Sheet="LTW"; Cell="N93"; Value=-4625.5}
    @{Sheet="LTW"; Cell="H127"; Value=49666.5}
    @{Sheet="LTW"; Cell="J127"; Value=49666.5}
    @{Sheet="LTW"; Cell="L127"; Value=49666.5}
    @{Sheet="LTW"; Cell="N127"; Value=-59586.5}
    @{Sheet="LTW"; Cell="H130"; Value=46147.332}
    @{Sheet="LTW"; Cell="J130"; Value=46147.332}
    @{Sheet="LTW"; Cell="L130"; Value=46147.332}
    @{Sheet="LTW"; Cell="N130"; Value=-56187.332}
    @{Sheet="LTW"; Cell="H139"; Value=36070.125}
    @{Sheet="LTW"; Cell="J139"; Value=36070.125}
    @{Sheet="LTW"; Cell="L139"; Value=36070.125}
    @{Sheet="LTW"; Cell="N139"; Value=-46350.125}
    @{Sheet="WVR"; Cell="H16"; Value=44750.5}
    @{Sheet="WVR"; Cell="J16"; Value=44750.5}
    @{Sheet="WVR"; Cell="L16"; Value=44750.5}
    @{Sheet="WVR"; Cell="N16"; Value=-45334.5}
    @{Sheet="WVR"; Cell="H137"; Value=53834.332}
    @{Sheet="WVR"; Cell="J137"; Value=53834.332}
    @{Sheet="WVR"; Cell="L137"; Value=53834.332}
    @{Sheet="WVR"; Cell="N137"; Value=-64034.332}
    @{Sheet="WVR"; Cell="H138"; Value=43936.43}
    @{Sheet="WVR"; Cell="J138"; Value=43936.43}
    @{Sheet="WVR"; Cell="L138"; Value=43936.43}
    @{Sheet="WVR"; Cell="N138"; Value=-54216.43}
    @{Sheet="WVR"; Cell="H139"; Value=34617.43}
    @{Sheet="WVR"; Cell="J139"; Value=34617.43}
    @{Sheet="WVR"; Cell="L139"; Value=34617.43}
    @{Sheet="WVR"; Cell="N139"; Value=-44897.43}
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}

# --- Cells that must be cleared entirely (no longer present after the edit) ---
$clears = @(
    @{Sheet="ARM"; Cell="N123"}
    @{Sheet="BSM"; Cell="N132"}
    @{Sheet="BSM"; Cell="N135"}
    @{Sheet="BSM"; Cell="N137"}
)

foreach ($c in $clears) {
    $ws = $wb.Worksheets.Item($c.Sheet)
    $ws.Range($c.Cell).ClearContents()
}
